$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: the numeric date serials (2022 fiscal demo dates) are
# replaced with literal text dates ("2023/5/15" etc). Switch the data
# cells to Text format first (numFmtId 49) so the strings are stored
# verbatim instead of being re-interpreted as numbers/dates, then write
# the new values. B2:B31 is handled before B1 so the plain "text" cell
# style is resolved/created ahead of the bordered header+text style,
# mirroring how the cells were actually edited.
$ws.Range("B2:B31").NumberFormat = "@"

$ws.Range("B2:B6").Value   = "2023/5/15"
$ws.Range("B7:B11").Value  = "2023/5/13"
$ws.Range("B12:B16").Value = "2023/5/12"
$ws.Range("B17:B21").Value = "2023/5/11"
$ws.Range("B22:B26").Value = "2023/5/10"
$ws.Range("B27:B31").Value = "2023/5/9"

# Header cell keeps its bordered/centered look but also becomes Text format.
$ws.Range("B1").NumberFormat = "@"

# --- Column widths ---
# ColumnWidth (character units) is rounded to the nearest pixel by the
# host, so feed it the value that rounds back to the desired stored width.
$ws.Columns.Item(2).ColumnWidth = 26.857142857142858   # -> stored width ~27.57 (closest to 27.625)
$ws.Columns.Item(4).ColumnWidth = 40.285714285714285   # -> stored width 41
$ws.Columns.Item(6).ColumnWidth = 20.428571427142857   # -> stored width ~21.14 (closest to 21.125)

# --- Selection ---
$ws.Range("B1").Select() | Out-Null
